$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '71.948.40'
$ws.Range("E2").Value = '  +4.79%  '
$ws.Range("D3").Value = '4.041.27'
$ws.Range("E3").Value = '  +4.69%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '539.51'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +3.22%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '153.35'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +8.71%  '
$ws.Range("E7").Value = '  +14.31%  '
$ws.Range("E8").Value = '  +0.04%  '
$ws.Range("E9").Value = '  +7.37%  '
$ws.Range("E10").Value = '  +4.22%  '
$ws.Range("E11").Value = '  +3.50%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '48.59'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +16.56%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '10.94'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +5.55%  '
$ws.Range("D14").Value = '4.684.89'
$ws.Range("E14").Value = '  +4.71%  '
$ws.Range("D15").Value = '4.055.81'
$ws.Range("E15").Value = '  +4.42%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.40'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +2.33%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '20.74'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -3.43%  '
$ws.Range("E18").Value = '  +1.74%  '
$ws.Range("E19").Value = '  -0.05%  '
$ws.Range("D20").Value = '71.953.94'
$ws.Range("E20").Value = '  +4.83%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '436.80'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +4.79%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '99.84'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +15.12%  '
$ws.Range("E23").Value = '  +2.08%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '4.28'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +6.29%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '14.71'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +4.70%  '
$ws.Range("E26").Value = '  -5.75%  '
$ws.Range("E27").Value = '  +5.23%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '3.70'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +30.01%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '37.20'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +4.89%  '
$ws.Range("E30").Value = '  +2.94%  '
$ws.Range("E31").Value = '  +2.51%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.133'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +6.46%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '682.08'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.95%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.94'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.47%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '67.24'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.41%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '43.12'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +9.24%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.435'
$ws.Range("D37").Style = "Normal"
$ws.Range("E38").Value = '  +7.19%  '
$ws.Range("D39").Value = '0.0₃0851'
$ws.Range("E39").Value = '  +0.01%  '
$ws.Range("E40").Value = '  +11.22%  '
$ws.Range("E41").Value = '  -1.39%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.999'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.07%  '
$ws.Range("E43").Value = '  +4.46%  '
$ws.Range("E44").Value = '  -0.08%  '
$ws.Range("E45").Value = '  +8.09%  '
$ws.Range("E46").Value = '  -4.78%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.44'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.75%  '
$ws.Range("E48").Value = '  +8.56%  '
$ws.Range("E49").Value = '  +2.50%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '3.39'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +3.31%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.000272'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -5.00%  '
